$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Range("E10").Value = 321
$ws.Range("E11").Value = 219
$ws.Range("E12").Value = 325
$ws.Range("E16").Value = 139
$ws.Range("E17").Value = 60
$ws.Range("E23").Value = 128
$ws.Range("E24").Value = 142
$ws.Range("E26").Value = 91
$ws.Range("E33").Value = 205
$ws.Range("E35").Value = 99
$ws.Range("E39").Value = 135
$ws.Range("E40").Value = 184
$ws.Range("E41").Value = 264
$ws.Range("E42").Value = 244
$ws.Range("E46").Value = 208
$ws.Range("E47").Value = 305
$ws.Range("E48").Value = 141
